$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H21").Value = 22800
$ws_ALC.Range("I21").Value = 20666.666
$ws_ALC.Range("J21").Value = 26000
$ws_ALC.Range("K21").Value = 20666.666
$ws_ALC.Range("L21").Value = 26000
$ws_ALC.Range("M21").Value = -20198.666
$ws_ALC.Range("N21").Value = -26936
$ws_ALC.Range("H23").Value = 22800
$ws_ALC.Range("I23").Value = 20666.666
$ws_ALC.Range("J23").Value = 26000
$ws_ALC.Range("K23").Value = 20666.666
$ws_ALC.Range("L23").Value = 26000
$ws_ALC.Range("M23").Value = -20432.666
$ws_ALC.Range("N23").Value = -26468
$ws_ALC.Range("H31").Value = 1276939.8
$ws_ALC.Range("I31").Value = 1276939.8
$ws_ALC.Range("J31").Value = 0
$ws_ALC.Range("K31").Value = 3830819.4
$ws_ALC.Range("L31").Value = 0
$ws_ALC.Range("M31").Value = -3830589.4
$ws_ALC.Range("N31").ClearContents()
$ws_ALC.Range("H38").Value = 785.2778
$ws_ALC.Range("I38").Value = 88.21429
$ws_ALC.Range("J38").Value = 3225
$ws_ALC.Range("K38").Value = 264.64287
$ws_ALC.Range("L38").Value = 9675
$ws_ALC.Range("M38").Value = 107.35713
$ws_ALC.Range("N38").Value = -10419
$ws_ALC.Range("H129").Value = 943.4375
$ws_ALC.Range("J129").Value = 1117.4783
$ws_ALC.Range("L129").Value = 3352.4349
$ws_ALC.Range("N129").Value = -13352.4349
$ws_ALC.Range("H132").Value = 2205.5356
$ws_ALC.Range("I132").Value = 1790.6444
$ws_ALC.Range("J132").Value = 3902.818
$ws_ALC.Range("K132").Value = 5371.933199999999
$ws_ALC.Range("L132").Value = 11708.454
$ws_ALC.Range("M132").Value = -2841.933199999999
$ws_ALC.Range("N132").Value = -16768.454
$ws_ALC.Range("H137").Value = 1093.5254
$ws_ALC.Range("I137").Value = 853.5111
$ws_ALC.Range("J137").Value = 1865
$ws_ALC.Range("K137").Value = 2560.5333
$ws_ALC.Range("L137").Value = 5595
$ws_ALC.Range("M137").Value = -10.53330000000005
$ws_ALC.Range("N137").Value = -10695
$ws_ALC.Range("H139").Value = 66180
$ws_ALC.Range("J139").Value = 66180
$ws_ALC.Range("L139").Value = 66180
$ws_ALC.Range("N139").Value = -76460
$ws_ARM.Range("H32").Value = 9228.789
$ws_ARM.Range("I32").Value = 9496.2
$ws_ARM.Range("K32").Value = 9496.2
$ws_ARM.Range("M32").Value = -9209.2
$ws_ARM.Range("H33").Value = 13001009
$ws_ARM.Range("I33").Value = 13001009
$ws_ARM.Range("K33").Value = 13001009
$ws_ARM.Range("M33").Value = -13000680
$ws_ARM.Range("H36").Value = 9996.667
$ws_ARM.Range("I36").Value = 9996.667
$ws_ARM.Range("K36").Value = 9996.667
$ws_ARM.Range("M36").Value = -9650.667
$ws_ARM.Range("H61").Value = 1035.7727
$ws_ARM.Range("I61").Value = 989.9048
$ws_ARM.Range("J61").Value = 1999
$ws_ARM.Range("K61").Value = 989.9048
$ws_ARM.Range("L61").Value = 1999
$ws_ARM.Range("M61").Value = -777.9048
$ws_ARM.Range("N61").Value = -2423
$ws_ARM.Range("H63").Value = 3347.5881
$ws_ARM.Range("I63").Value = 2931.8125
$ws_ARM.Range("K63").Value = 2931.8125
$ws_ARM.Range("M63").Value = -2245.8125
$ws_ARM.Range("H66").Value = 3347.5881
$ws_ARM.Range("I66").Value = 2931.8125
$ws_ARM.Range("K66").Value = 14659.0625
$ws_ARM.Range("M66").Value = -11227.0625
$ws_ARM.Range("H74").Value = 1131.069
$ws_ARM.Range("I74").Value = 912.85
$ws_ARM.Range("J74").Value = 1616
$ws_ARM.Range("K74").Value = 912.85
$ws_ARM.Range("L74").Value = 1616
$ws_ARM.Range("M74").Value = -38.85000000000002
$ws_ARM.Range("N74").Value = -3364
$ws_ARM.Range("H77").Value = 1131.069
$ws_ARM.Range("I77").Value = 912.85
$ws_ARM.Range("J77").Value = 1616
$ws_ARM.Range("K77").Value = 4564.25
$ws_ARM.Range("L77").Value = 8080
$ws_ARM.Range("M77").Value = -196.25
$ws_ARM.Range("N77").Value = -16816
$ws_ARM.Range("H132").Value = 870463.6
$ws_ARM.Range("I132").Value = 909802.94
$ws_ARM.Range("K132").Value = 2729408.82
$ws_ARM.Range("M132").Value = -2726878.82
$ws_ARM.Range("H136").Value = 1035.7727
$ws_ARM.Range("I136").Value = 989.9048
$ws_ARM.Range("J136").Value = 1999
$ws_ARM.Range("K136").Value = 2969.7144
$ws_ARM.Range("L136").Value = 5997
$ws_ARM.Range("M136").Value = -419.7143999999998
$ws_ARM.Range("N136").Value = -11097
$ws_BSM.Range("H134").Value = 608787.94
$ws_BSM.Range("I134").Value = 771461.56
$ws_BSM.Range("J134").Value = 4571.4287
$ws_BSM.Range("K134").Value = 2314384.68
$ws_BSM.Range("L134").Value = 13714.2861
$ws_BSM.Range("M134").Value = -2311849.68
$ws_BSM.Range("N134").Value = -18784.2861
$ws_BSM.Range("H140").Value = 53032.555
$ws_BSM.Range("J140").Value = 53032.555
$ws_BSM.Range("L140").Value = 53032.555
$ws_BSM.Range("N140").Value = -63392.555
$ws_CRP.Range("H29").Value = 50000
$ws_CRP.Range("J29").Value = 0
$ws_CRP.Range("L29").Value = 0
$ws_CRP.Range("N29").ClearContents()
$ws_CRP.Range("H31").Value = 17243744
$ws_CRP.Range("I31").Value = 22728946
$ws_CRP.Range("J31").Value = 4541.2856
$ws_CRP.Range("K31").Value = 22728946
$ws_CRP.Range("L31").Value = 4541.2856
$ws_CRP.Range("M31").Value = -22728651
$ws_CRP.Range("N31").Value = -5131.2856
$ws_CRP.Range("H34").Value = 17243744
$ws_CRP.Range("I34").Value = 22728946
$ws_CRP.Range("J34").Value = 4541.2856
$ws_CRP.Range("K34").Value = 22728946
$ws_CRP.Range("L34").Value = 4541.2856
$ws_CRP.Range("M34").Value = -22728744
$ws_CRP.Range("N34").Value = -4945.2856
$ws_CRP.Range("H58").Value = 1230.359
$ws_CRP.Range("I58").Value = 1182.1666
$ws_CRP.Range("K58").Value = 1182.1666
$ws_CRP.Range("M58").Value = -979.1666
$ws_CRP.Range("H107").Value = 648.4286
$ws_CRP.Range("I107").Value = 609.75
$ws_CRP.Range("K107").Value = 609.75
$ws_CRP.Range("M107").Value = 1310.25
$ws_CRP.Range("H132").Value = 2226.0293
$ws_CRP.Range("I132").Value = 1851.6666
$ws_CRP.Range("K132").Value = 5554.9998
$ws_CRP.Range("M132").Value = -3024.9998
$ws_CRP.Range("H136").Value = 1230.359
$ws_CRP.Range("I136").Value = 1182.1666
$ws_CRP.Range("K136").Value = 3546.4998
$ws_CRP.Range("M136").Value = -996.4998
$ws_CUL.Range("H41").Value = 875.68
$ws_CUL.Range("J41").Value = 895.8333
$ws_CUL.Range("L41").Value = 2687.4999
$ws_CUL.Range("N41").Value = -3363.4999
$ws_GSM.Range("H132").Value = 2096.3076
$ws_GSM.Range("I132").Value = 1694.5161
$ws_GSM.Range("J132").Value = 3653.25
$ws_GSM.Range("K132").Value = 5083.5483
$ws_GSM.Range("L132").Value = 10959.75
$ws_GSM.Range("M132").Value = -2553.5483
$ws_GSM.Range("N132").Value = -16019.75
$ws_LTW.Range("H34").Value = 49663.332
$ws_LTW.Range("I34").Value = 49663.332
$ws_LTW.Range("K34").Value = 49663.332
$ws_LTW.Range("M34").Value = -49491.332
$ws_LTW.Range("H132").Value = 5498.4644
$ws_LTW.Range("I132").Value = 5769.913
$ws_LTW.Range("K132").Value = 17309.739
$ws_LTW.Range("M132").Value = -14779.739
$ws_LTW.Range("H136").Value = 1858.725
$ws_LTW.Range("I136").Value = 1873.4062
$ws_LTW.Range("J136").Value = 1800
$ws_LTW.Range("K136").Value = 5620.2186
$ws_LTW.Range("L136").Value = 5400
$ws_LTW.Range("M136").Value = -3070.2186
$ws_LTW.Range("N136").Value = -10500
$ws_WVR.Range("H26").Value = 50367.6
$ws_WVR.Range("I26").Value = 37956
$ws_WVR.Range("J26").Value = 100014
$ws_WVR.Range("K26").Value = 37956
$ws_WVR.Range("L26").Value = 100014
$ws_WVR.Range("M26").Value = -37663
$ws_WVR.Range("N26").Value = -100600
$ws_WVR.Range("H43").Value = 49999.5
$ws_WVR.Range("I43").Value = 49999.5
$ws_WVR.Range("K43").Value = 49999.5
$ws_WVR.Range("M43").Value = -49850.5
$ws_WVR.Range("H132").Value = 2354.25
$ws_WVR.Range("I132").Value = 1763.3529
$ws_WVR.Range("J132").Value = 3023.9333
$ws_WVR.Range("K132").Value = 5290.0587
$ws_WVR.Range("L132").Value = 9071.7999
$ws_WVR.Range("M132").Value = -2760.0587
$ws_WVR.Range("N132").Value = -14131.7999
$ws_WVR.Range("H136").Value = 1291.5526
$ws_WVR.Range("I136").Value = 1211.12
$ws_WVR.Range("J136").Value = 1446.2307
$ws_WVR.Range("K136").Value = 3633.36
$ws_WVR.Range("L136").Value = 4338.6921
$ws_WVR.Range("M136").Value = -1083.36
$ws_WVR.Range("N136").Value = -9438.6921
$ws_WVR.Range("H141").Value = 75459
$ws_WVR.Range("J141").Value = 75459
$ws_WVR.Range("L141").Value = 75459
$ws_WVR.Range("N141").Value = -85819
